$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the columns/rows that are no longer part of the (now smaller) table
$ws.Range("D1:E7").Clear()
$ws.Range("A4:C7").Clear()

# Update the remaining regression coefficients to their new values
$ws.Range("B2").Value = "-0.865***"
$ws.Range("C2").Value = "0.357**"
$ws.Range("B3").Value = "-0.256***"

# "0.009" looks numeric, so force it to be stored as text (matching the
# other coefficient cells), then strip the number-format override so the
# cell keeps the default (un-styled) appearance
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.009"
$ws.Range("C3").Style = "Normal"
